$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 31665
$ws.Range("K2").Value = 38004
$ws.Range("L2").Value = 652
$ws.Range("M2").Value = 870
$ws.Range("J3").Value = 15950
$ws.Range("K3").Value = 21772
$ws.Range("L3").Value = 513
$ws.Range("M3").Value = 116
$ws.Range("J4").Value = 44678
$ws.Range("K4").Value = 73561
$ws.Range("L4").Value = 2918
$ws.Range("M4").Value = 297
$ws.Range("J5").Value = 1340
$ws.Range("K5").Value = 1392
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 14
$ws.Range("J6").Value = 32919
$ws.Range("K6").Value = 44048
$ws.Range("L6").Value = 1013
$ws.Range("M6").Value = 683
$ws.Range("J7").Value = 4080
$ws.Range("K7").Value = 3597
$ws.Range("L7").Value = 529
$ws.Range("M7").Value = 29
$ws.Range("J8").Value = 5872
$ws.Range("K8").Value = 3500
$ws.Range("L8").Value = 43
$ws.Range("M8").Value = 0
$ws.Range("J9").Value = 2651
$ws.Range("K9").Value = 1866
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("J10").Value = 144
$ws.Range("K10").Value = 310
$ws.Range("M10").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 2
$ws.Range("J12").Value = 790
$ws.Range("K12").Value = 575
$ws.Range("L12").Value = 71
$ws.Range("M12").Value = 244
$ws.Range("J13").Value = 2985
$ws.Range("K13").Value = 1679
$ws.Range("L13").Value = 272
$ws.Range("M13").Value = 195
$ws.Range("J14").Value = 4439
$ws.Range("K14").Value = 3365
$ws.Range("L14").Value = 338
$ws.Range("M14").Value = 43
$ws.Range("J15").Value = 2764
$ws.Range("K15").Value = 2931
$ws.Range("L15").Value = 154
$ws.Range("M15").Value = 34
$ws.Range("J16").Value = 1556
$ws.Range("K16").Value = 1950
$ws.Range("L16").Value = 65
$ws.Range("M16").Value = 0
$ws.Range("J17").Value = 10200
$ws.Range("K17").Value = 17902
$ws.Range("L17").Value = 141
$ws.Range("M17").Value = 144
$ws.Range("J18").Value = 721
$ws.Range("K18").Value = 1704
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("J19").Value = 13766
$ws.Range("K19").Value = 16100
$ws.Range("L19").Value = 188
$ws.Range("M19").Value = 0
$ws.Range("J20").Value = 339
$ws.Range("K20").Value = 46
$ws.Range("L20").Value = 95
$ws.Range("M20").Value = 0
$ws.Range("J21").Value = 16184
$ws.Range("K21").Value = 11360
$ws.Range("L21").Value = 419
$ws.Range("M21").Value = 0
$ws.Range("J22").Value = 1363
$ws.Range("K22").Value = 461
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = 0
$ws.Range("J23").Value = 13504
$ws.Range("K23").Value = 18309
$ws.Range("L23").Value = 842
$ws.Range("M23").Value = 19
$ws.Range("J24").Value = 59749
$ws.Range("K24").Value = 65286
$ws.Range("L24").Value = 2910
$ws.Range("M24").Value = 676
$ws.Range("J25").Value = 3697
$ws.Range("K25").Value = 5871
$ws.Range("L25").Value = 477
$ws.Range("M25").Value = 0
$ws.Range("J27").Value = 5199
$ws.Range("K27").Value = 3137
$ws.Range("L27").Value = 385
$ws.Range("M27").Value = 1
$ws.Range("J28").Value = 943
$ws.Range("K28").Value = 1361
$ws.Range("L28").Value = 18
$ws.Range("M28").Value = 0
$ws.Range("J29").Value = 9306
$ws.Range("K29").Value = 14164
$ws.Range("L29").Value = 420
$ws.Range("M29").Value = 0
$ws.Range("J30").Value = 483
$ws.Range("K30").Value = 215
$ws.Range("L30").Value = 18
$ws.Range("M30").Value = 16
$ws.Range("J31").Value = 1592
$ws.Range("K31").Value = 2628
$ws.Range("L31").Value = 201
$ws.Range("M31").Value = 0
$ws.Range("J32").Value = 14772
$ws.Range("K32").Value = 8492
$ws.Range("L32").Value = 2751
$ws.Range("M32").Value = 306
$ws.Range("J33").Value = 8931
$ws.Range("K33").Value = 7673
$ws.Range("L33").Value = 129
$ws.Range("M33").Value = 0
$ws.Range("J34").Value = 4508
$ws.Range("K34").Value = 3164
$ws.Range("L34").Value = 1573
$ws.Range("M34").Value = 0
$ws.Range("J35").Value = 36080
$ws.Range("K35").Value = 54815
$ws.Range("L35").Value = 1753
$ws.Range("M35").Value = 219
$ws.Range("J36").Value = 8100
$ws.Range("K36").Value = 5435
$ws.Range("L36").Value = 186
$ws.Range("M36").Value = 77
$ws.Range("J37").Value = 14659
$ws.Range("K37").Value = 22937
$ws.Range("L37").Value = 4579
$ws.Range("M37").Value = 0
$ws.Range("J38").Value = 569
$ws.Range("K38").Value = 991
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 4
$ws.Range("J39").Value = 434
$ws.Range("K39").Value = 2625
$ws.Range("L39").Value = 38
$ws.Range("M39").Value = 1
$ws.Range("J40").Value = 135
$ws.Range("K40").Value = 4619
$ws.Range("L40").Value = 5
$ws.Range("M40").Value = 0
$ws.Range("J41").Value = 98
$ws.Range("K41").Value = 17642
$ws.Range("L41").Value = 41
$ws.Range("M41").Value = 1
$ws.Range("J42").Value = 379
$ws.Range("K42").Value = 44
$ws.Range("L42").Value = 3
$ws.Range("M42").Value = 0
$ws.Range("J43").Value = 35
$ws.Range("K43").Value = 1406
$ws.Range("M43").Value = 0
$ws.Range("J44").Value = 41
$ws.Range("K44").Value = 998
$ws.Range("L44").Value = 6
$ws.Range("M44").Value = 0
$ws.Range("J45").Value = 2098
$ws.Range("K45").Value = 1379
$ws.Range("L45").Value = 2044
$ws.Range("M45").Value = 2
$ws.Range("J46").Value = 9958
$ws.Range("K46").Value = 10376
$ws.Range("L46").Value = 44
$ws.Range("M46").Value = 458
$ws.Range("J47").Value = 12387
$ws.Range("K47").Value = 36684
$ws.Range("L47").Value = 1021
$ws.Range("M47").Value = 296
$ws.Range("J48").Value = 9542
$ws.Range("K48").Value = 14093
$ws.Range("L48").Value = 87
$ws.Range("M48").Value = 0
$ws.Range("J49").Value = 5887
$ws.Range("K49").Value = 12024
$ws.Range("L49").Value = 552
$ws.Range("M49").Value = 165
$ws.Range("J50").Value = 22025
$ws.Range("K50").Value = 25535
$ws.Range("L50").Value = 3743
$ws.Range("M50").Value = 708
$ws.Range("J51").Value = 4080
$ws.Range("K51").Value = 3779
$ws.Range("L51").Value = 9
$ws.Range("M51").Value = 0
$ws.Range("J52").Value = 8248
$ws.Range("K52").Value = 13044
$ws.Range("L52").Value = 656
$ws.Range("M52").Value = 2
$ws.Range("J53").Value = 2214
$ws.Range("K53").Value = 1058
$ws.Range("L53").Value = 37
$ws.Range("M53").Value = 0
$ws.Range("J54").Value = 1510
$ws.Range("K54").Value = 2094
$ws.Range("L54").Value = 40
$ws.Range("M54").Value = 0
$ws.Range("J55").Value = 2874
$ws.Range("K55").Value = 4032
$ws.Range("L55").Value = 5
$ws.Range("M55").Value = 3
